$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Disease Ontology row: update source_version from v2025-03-31 to v2025-04-29
$ws.Range("E3").Value = "v2025-04-29"

# Update the active cell selection to E3 (matches the captured UI state)
$ws.Range("E3").Select()
